# Daily attendance processing - 2025-12-12 13:43:33
# Normalize the "Recorded By" (column G) entries so that "System" is
# listed in the position the other identity (email) previously occupied
# among the first two comma-separated tokens - i.e. swap the first two
# tokens whenever "System" is one of them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $raw = $cell.Value2
    if ($null -eq $raw) { continue }

    $value = [string]$raw
    if ($value -eq "") { continue }

    $parts = $value -split ", "
    if ($parts.Count -ge 2 -and ($parts[0] -eq "System" -or $parts[1] -eq "System")) {
        $tmp = $parts[0]
        $parts[0] = $parts[1]
        $parts[1] = $tmp
        $newValue = [string]::Join(", ", $parts)
        $cell.Value = $newValue
    }
}
